# Refresh the "cryptos" price list (GitHub Actions scheduled update).
# For numeric-looking Price values we prefix with a leading apostrophe so
# Excel keeps them as literal text (preserving trailing zeros / thousands
# dots) instead of silently coercing them to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '25.939.62'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '1.637.23'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D5').Value = "'215.49"
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = "'19.59"
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('D11').Value = "'0.0794"
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '1.866.14'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '1.649.75'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = "'62.87"
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '25.919.27'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').Value = "'192.63"
$ws.Range('E20').Value = '  -1.27%  '
$ws.Range('D21').Value = "'4.35"
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('D23').Value = "'6.26"
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  +5.09%  '
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'143.21"
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').Value = "'6.87"
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('D29').Value = "'15.56"
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').Value = "'0.0500"
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').Value = "'3.28"
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  -4.92%  '
$ws.Range('E35').Value = '  +1.81%  '
$ws.Range('D36').Value = "'0.900"
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('D37').Value = '1.131.89'
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').Value = "'5.47"
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').Value = "'99.22"
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').Value = "'0.796"
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').Value = '1.775.68'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('D46').Value = "'56.59"
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').Value = "'0.0529"
$ws.Range('E47').Value = '  +2.18%  '
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').Value = "'7.68"
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').Value = "'0.414"
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').Value = "'0.0959"
$ws.Range('E51').Value = '  -1.01%  '
